# Aufwand_git.xlsx - add a new time-tracking entry row (row 67) on the
# "Tabelle1" sheet, mirroring the existing "Performance issues" task from
# row 66 (the bugfix text was generated even though the shape wasn't
# actually included, per the commit message).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New entry: 2024-04-02, 7 hours, same task/comment text as row 66.
$ws.Range("A67").Value = 45384
$ws.Range("A67").NumberFormat = $ws.Range("A66").NumberFormat
$ws.Range("B67").Value = 7
$ws.Range("C67").Value = $ws.Range("C66").Value2

# Move the active selection past the newly added row, like Excel does
# after typing data into the last row of a contiguous block.
$ws.Activate()
$ws.Range("A68").Select()
